$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data set: 6 rows (r2..r7) spanning columns A..T
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol, D=Target cluster,
# E..T = numeric metrics

$rows = @(
  @{ Row=2; A="ECs";  D="ECs";  E=3; F=1; G=297.8183156666666;  H=893.4549469999999;  I=0.8852156413092672;  J=0.8852156413092673;  K=3; L=1;                  M=7.879565666666667; N=23.638697; O=0.9977172793687663;  P=0.9977172793687664;  Q=2346.678975031562;  R=21120.11077528406;  S=0.8831949413017598;    T=0.88319494130176 },
  @{ Row=3; A="ECs";  D="sCs";  E=3; F=1; G=297.8183156666666;  H=893.4549469999999;  I=0.8852156413092672;  J=0.8852156413092673;  K=1; L=0.3333333333333333; M=0.018028;          N=0.054084;  O=0.002282720631233623; P=0.002282720631233623; Q=5.369068594838666;  R=48.321617353548;    S=0.002020700007507367;  T=0.002020700007507367 },
  @{ Row=4; A="FAPs"; D="ECs";  E=3; F=1; G=24.34034433333333;  H=73.021033;          I=0.07234764413494278; J=0.0723476441349428;  K=3; L=1;                  M=7.879565666666667; N=23.638697; O=0.9977172793687663;  P=0.9977172793687664;  Q=191.7913415237779;  R=1726.122073714001;  S=0.07218249467505479;   T=0.07218249467505482 },
  @{ Row=5; A="FAPs"; D="sCs";  E=3; F=1; G=24.34034433333333;  H=73.021033;          I=0.07234764413494278; J=0.0723476441349428;  K=1; L=0.3333333333333333; M=0.018028;          N=0.054084;  O=0.002282720631233623; P=0.002282720631233623; Q=0.4388077276413333; R=3.949269548772;     S=0.0001651494598879821; T=0.0001651494598879822 },
  @{ Row=6; A="sCs";  D="ECs";  E=3; F=1; G=14.277234;          H=42.831702;          I=0.04243671455578994; J=0.04243671455578994; K=3; L=1;                  M=7.879565666666667; N=23.638697; O=0.9977172793687663;  P=0.9977172793687664;  Q=112.498402841366;   R=1012.485625572294;  S=0.04233984339195166;   T=0.04233984339195167 },
  @{ Row=7; A="sCs";  D="sCs";  E=3; F=1; G=14.277234;          H=42.831702;          I=0.04243671455578994; J=0.04243671455578994; K=1; L=0.3333333333333333; M=0.018028;          N=0.054084;  O=0.002282720631233623; P=0.002282720631233623; Q=0.257389974552;     R=2.316509770968;     S=0.0000968711638382739; T=0.00009687116383827391 }
)

foreach ($r in $rows) {
  $row = $r.Row
  $ws.Cells.Item($row, 1).Value = $r.A
  $ws.Cells.Item($row, 2).Value = "Cd44"
  $ws.Cells.Item($row, 3).Value = "Sele"
  $ws.Cells.Item($row, 4).Value = $r.D
  $ws.Cells.Item($row, 5).Value = $r.E
  $ws.Cells.Item($row, 6).Value = $r.F
  $ws.Cells.Item($row, 7).Value = $r.G
  $ws.Cells.Item($row, 8).Value = $r.H
  $ws.Cells.Item($row, 9).Value = $r.I
  $ws.Cells.Item($row, 10).Value = $r.J
  $ws.Cells.Item($row, 11).Value = $r.K
  $ws.Cells.Item($row, 12).Value = $r.L
  $ws.Cells.Item($row, 13).Value = $r.M
  $ws.Cells.Item($row, 14).Value = $r.N
  $ws.Cells.Item($row, 15).Value = $r.O
  $ws.Cells.Item($row, 16).Value = $r.P
  $ws.Cells.Item($row, 17).Value = $r.Q
  $ws.Cells.Item($row, 18).Value = $r.R
  $ws.Cells.Item($row, 19).Value = $r.S
  $ws.Cells.Item($row, 20).Value = $r.T
}
